# Tool Guide update: append " and asking questions." right after
# "... people start testing" (and before the closing parenthesis) in the
# "QUEUE FOR NEXT DAY - FAQs" intro paragraph.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "FAQs will keep on getting updated once people start testing", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the sentence to extend."
}

# Collapse to the end of the matched text (just before the closing ")")
# and insert the new sentence fragment there, inheriting the surrounding
# (identical) run formatting.
$rng.Collapse(0)
$rng.InsertAfter(" and asking questions.")
